$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 15, shifting rows 15:66 down to 16:67
$ws.Rows.Item(15).Insert()

# Fill in the new row 15 with the new data point
$ws.Cells.Item(15, 1).Value = 11
$ws.Cells.Item(15, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(15, 3).Value = "Bíobío"
$ws.Cells.Item(15, 4).Value = 44453
$ws.Cells.Item(15, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(15, 5).Value = 8
$ws.Cells.Item(15, 6).Value = 100112043
$ws.Cells.Item(15, 7).Value = "Pepino ensalada"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 100
$ws.Cells.Item(15, 11).Value = 16000
$ws.Cells.Item(15, 12).Value = 17000
$ws.Cells.Item(15, 13).Value = 16500
$ws.Cells.Item(15, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(15, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(15, 16).Value = 275
$ws.Cells.Item(15, 17).Value = 60
$ws.Cells.Item(15, 18).Value = "Hortaliza"
